$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (G=5489)
$ws.Range("H2").Value = 949.3333
$ws.Range("I2").Value = 50
$ws.Range("K2").Value = 50
$ws.Range("M2").Value = 63

# Row 28 (G=27772)
$ws.Range("H28").Value = 565.45
$ws.Range("I28").Value = 430
$ws.Range("J28").Value = 1333
$ws.Range("K28").Value = 430
$ws.Range("L28").Value = 1333
$ws.Range("M28").Value = 55
$ws.Range("N28").Value = -2303

# Row 32 (G=5484)
$ws.Range("H32").Value = 3924.1765
$ws.Range("I32").Value = 3700.1
$ws.Range("J32").Value = 4244.2856
$ws.Range("K32").Value = 3700.1
$ws.Range("L32").Value = 4244.2856
$ws.Range("M32").Value = -3374.1
$ws.Range("N32").Value = -4896.2856

# Row 41 (G=5478)
$ws.Range("H41").Value = 52992
$ws.Range("I41").Value = 148
$ws.Range("J41").Value = 125652.5
$ws.Range("K41").Value = 148
$ws.Range("L41").Value = 125652.5
$ws.Range("M41").Value = 292
$ws.Range("N41").Value = -126532.5

# Row 51 (G=5486)
$ws.Range("H51").Value = 27521.875
$ws.Range("I51").Value = 52287
$ws.Range("J51").Value = 19266.834
$ws.Range("K51").Value = 52287
$ws.Range("L51").Value = 19266.834
$ws.Range("M51").Value = -51803
$ws.Range("N51").Value = -20234.834

# Row 61 (G=4604)
$ws.Range("H61").Value = 2499
$ws.Range("J61").Value = 2999
$ws.Range("L61").Value = 8997
$ws.Range("N61").Value = -9341

# Row 98 (G=36237)
$ws.Range("H98").Value = 1466.3414
$ws.Range("I98").Value = 1510.4736
$ws.Range("K98").Value = 1510.4736
$ws.Range("M98").Value = -12.47360000000003

# Row 112 (G=27960)
$ws.Range("H112").Value = 1835.4546
$ws.Range("I112").Value = 1269.4
$ws.Range("J112").Value = 2001.9412
$ws.Range("K112").Value = 3808.2
$ws.Range("L112").Value = 6005.8236
$ws.Range("M112").Value = -2700.2
$ws.Range("N112").Value = -8221.8236

# Row 122 (G=36237)
$ws.Range("H122").Value = 1466.3414
$ws.Range("I122").Value = 1510.4736
$ws.Range("K122").Value = 4531.4208
$ws.Range("M122").Value = -2081.4208

# Row 137 (G=44013)
$ws.Range("H137").Value = 3137.1562
$ws.Range("J137").Value = 3905.6
$ws.Range("L137").Value = 11716.8
$ws.Range("N137").Value = -16816.8

# Row 138 (G=44169)
$ws.Range("H138").Value = 3637.0156
$ws.Range("I138").Value = 3054.9333
$ws.Range("J138").Value = 3815.204
$ws.Range("K138").Value = 9164.7999
$ws.Range("L138").Value = 11445.612
$ws.Range("M138").Value = -4024.7999
$ws.Range("N138").Value = -21725.612

# Row 141 (G=44161)
$ws.Range("H141").Value = 1304.6
$ws.Range("I141").Value = 1227.3334
$ws.Range("K141").Value = 3682.0002
$ws.Range("M141").Value = 1497.9998

$ws = $wb.Worksheets.Item("ARM")
# Row 4 (G=5071)
$ws.Range("H4").Value = 676.55554
$ws.Range("I4").Value = 165.5
$ws.Range("K4").Value = 165.5
$ws.Range("M4").Value = -49.5

# Row 10 (G=2662)
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

# Row 45 (G=27714)
$ws.Range("H45").Value = 4375.25
$ws.Range("I45").Value = 4240.3
$ws.Range("K45").Value = 4240.3
$ws.Range("M45").Value = -3863.3

# Row 63 (G=12528)
$ws.Range("H63").Value = 6573
$ws.Range("I63").Value = 2909.2
$ws.Range("J63").Value = 8862.875
$ws.Range("K63").Value = 2909.2
$ws.Range("L63").Value = 8862.875
$ws.Range("M63").Value = -2223.2
$ws.Range("N63").Value = -10234.875

# Row 66 (G=12528)
$ws.Range("H66").Value = 6573
$ws.Range("I66").Value = 2909.2
$ws.Range("J66").Value = 8862.875
$ws.Range("K66").Value = 14546
$ws.Range("L66").Value = 44314.375
$ws.Range("M66").Value = -11114
$ws.Range("N66").Value = -51178.375

# Row 102 (G=19945)
$ws.Range("H102").Value = 6358.0835
$ws.Range("I102").Value = 6629.8
$ws.Range("J102").Value = 4999.5
$ws.Range("K102").Value = 6629.8
$ws.Range("L102").Value = 4999.5
$ws.Range("M102").Value = -5007.8
$ws.Range("N102").Value = -8243.5

# Row 132 (G=43997)
$ws.Range("H132").Value = 634814.5
$ws.Range("I132").Value = 898721.0600000001
$ws.Range("J132").Value = 11035.363
$ws.Range("K132").Value = 2696163.18
$ws.Range("L132").Value = 33106.089
$ws.Range("M132").Value = -2693633.18
$ws.Range("N132").Value = -38166.089

$ws = $wb.Worksheets.Item("CRP")
# Row 50 (G=1862)
$ws.Range("H50").Value = 119996
$ws.Range("J50").Value = 119996
$ws.Range("L50").Value = 119996
$ws.Range("N50").Value = -121246

# Row 51 (G=2039)
$ws.Range("H51").Value = 130063
$ws.Range("J51").Value = 130063
$ws.Range("L51").Value = 130063
$ws.Range("N51").Value = -131535

# Row 60 (G=1937)
$ws.Range("H60").Value = 78481.45
$ws.Range("J60").Value = 82379.60000000001
$ws.Range("L60").Value = 82379.60000000001
$ws.Range("N60").Value = -83401.60000000001

# Row 61 (G=2039)
$ws.Range("H61").Value = 130063
$ws.Range("J61").Value = 130063
$ws.Range("L61").Value = 130063
$ws.Range("N61").Value = -130759

# Row 62 (G=12580)
$ws.Range("H62").Value = 4576.0713
$ws.Range("I62").Value = 4343.222
$ws.Range("J62").Value = 4995.2
$ws.Range("K62").Value = 4343.222
$ws.Range("L62").Value = 4995.2
$ws.Range("M62").Value = -3719.222
$ws.Range("N62").Value = -6243.2

# Row 63 (G=10604)
$ws.Range("H63").Value = 77847.164
$ws.Range("J63").Value = 78816.8
$ws.Range("L63").Value = 78816.8
$ws.Range("N63").Value = -80188.8

# Row 65 (G=12580)
$ws.Range("H65").Value = 4576.0713
$ws.Range("I65").Value = 4343.222
$ws.Range("J65").Value = 4995.2
$ws.Range("K65").Value = 21716.11
$ws.Range("L65").Value = 24976
$ws.Range("M65").Value = -18596.11
$ws.Range("N65").Value = -31216

# Row 66 (G=10604)
$ws.Range("H66").Value = 77847.164
$ws.Range("J66").Value = 78816.8
$ws.Range("L66").Value = 236450.4
$ws.Range("N66").Value = -243314.4

# Row 124 (G=34285)
$ws.Range("H124").Value = 77663
$ws.Range("J124").Value = 77663
$ws.Range("L124").Value = 77663
$ws.Range("N124").Value = -82573

# Row 132 (G=44019)
$ws.Range("H132").Value = 7589699.5
$ws.Range("I132").Value = 16548.416
$ws.Range("K132").Value = 49645.24800000001
$ws.Range("M132").Value = -47115.24800000001

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (G=36249)
$ws.Range("H7").Value = 4168.857
$ws.Range("I7").Value = 3924
$ws.Range("J7").Value = 6495
$ws.Range("K7").Value = 3924
$ws.Range("L7").Value = 6495
$ws.Range("M7").Value = -3812
$ws.Range("N7").Value = -6719

# Row 55 (G=5284)
$ws.Range("H55").Value = 1898.1818
$ws.Range("I55").Value = 820
$ws.Range("J55").Value = 2514.2856
$ws.Range("K55").Value = 820
$ws.Range("L55").Value = 2514.2856
$ws.Range("M55").Value = -647
$ws.Range("N55").Value = -2860.2856

# Row 61 (G=27740)
$ws.Range("H61").Value = 2322.524
$ws.Range("J61").Value = 3225.75
$ws.Range("L61").Value = 3225.75
$ws.Range("N61").Value = -3629.75

# Row 68 (G=12563)
$ws.Range("H68").Value = 6862.6665
$ws.Range("I68").Value = 6425.4287
$ws.Range("J68").Value = 7474.8
$ws.Range("K68").Value = 6425.4287
$ws.Range("L68").Value = 7474.8
$ws.Range("M68").Value = -5676.4287
$ws.Range("N68").Value = -8972.799999999999

# Row 71 (G=12563)
$ws.Range("H71").Value = 6862.6665
$ws.Range("I71").Value = 6425.4287
$ws.Range("J71").Value = 7474.8
$ws.Range("K71").Value = 32127.1435
$ws.Range("L71").Value = 37374
$ws.Range("M71").Value = -28383.1435
$ws.Range("N71").Value = -44862

# Row 82 (G=12565)
$ws.Range("H82").Value = 5524.3
$ws.Range("I82").Value = 3573.8333
$ws.Range("K82").Value = 3573.8333
$ws.Range("M82").Value = -3212.8333

# Row 85 (G=12565)
$ws.Range("H85").Value = 5524.3
$ws.Range("I85").Value = 3573.8333
$ws.Range("K85").Value = 3573.8333
$ws.Range("M85").Value = -2325.8333

# Row 113 (G=27740)
$ws.Range("H113").Value = 2322.524
$ws.Range("J113").Value = 3225.75
$ws.Range("L113").Value = 3225.75
$ws.Range("N113").Value = -7565.75

# Row 126 (G=36249)
$ws.Range("H126").Value = 4168.857
$ws.Range("I126").Value = 3924
$ws.Range("J126").Value = 6495
$ws.Range("K126").Value = 11772
$ws.Range("L126").Value = 19485
$ws.Range("M126").Value = -9302
$ws.Range("N126").Value = -24425

# Row 132 (G=44058)
$ws.Range("H132").Value = 1120595.1
$ws.Range("I132").Value = 1575975.2
$ws.Range("J132").Value = 7443.778
$ws.Range("K132").Value = 4727925.6
$ws.Range("L132").Value = 22331.334
$ws.Range("M132").Value = -4725395.6
$ws.Range("N132").Value = -27391.334

# Row 136 (G=44060)
$ws.Range("H136").Value = 9166.666999999999
$ws.Range("I136").Value = 8000
$ws.Range("K136").Value = 24000
$ws.Range("M136").Value = -21450

$ws = $wb.Worksheets.Item("WVR")
# Row 113 (G=27752)
$ws.Range("H113").Value = 1731.6522
$ws.Range("I113").Value = 1125.2307
$ws.Range("J113").Value = 2520
$ws.Range("K113").Value = 3375.6921
$ws.Range("L113").Value = 7560
$ws.Range("M113").Value = -1205.6921
$ws.Range("N113").Value = -11900

# Row 126 (G=36210)
$ws.Range("H126").Value = 3660.3333
$ws.Range("I126").Value = 3592.2942
$ws.Range("K126").Value = 10776.8826
$ws.Range("M126").Value = -8306.882599999999
